# Update Name of Algo
# Applies updated imputed values in columns A and D for the KNN result data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"  = -20.27
    "D4"  = -7.786999999999999
    "A6"  = -22.311
    "A7"  = -20.115
    "D9"  = -7.797
    "D12" = -7.072
    "A16" = -22.107
    "D17" = -8.431999999999999
    "D18" = -8.653
    "D19" = -8.047999999999998
    "A20" = -20.159
    "D20" = -7.819999999999999
    "D26" = -7.628
    "A28" = -21.927
    "A29" = -21.439
    "D31" = -7.984
    "A32" = -21.713
    "D39" = -7.699
    "A40" = -19.841
    "D40" = -8.16
    "D41" = -7.928
    "D42" = -7.958
    "D43" = -7.877999999999998
    "A46" = -21.801
    "D47" = -7.569000000000001
    "D48" = -7.717999999999999
    "A51" = -21.95
    "A52" = -22.036
    "A57" = -22.339
    "A59" = -22.699
    "A62" = -22.087
    "D63" = -7.178999999999999
    "D64" = -7.068
    "A66" = -21.615
    "A73" = -20.213
    "A74" = -21.246
    "D76" = -7.891000000000001
    "D81" = -7.674000000000001
    "D89" = -8.172000000000001
    "A92" = -21.566
    "D94" = -7.515000000000001
    "A100" = -22.465
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
